# Wrote custom step aic (stepAIC_custom.R) for general models
#
# The workbook has columns: year, NTS, NTE, NTW, KL, HK, NTN (A:G)
# The "NTE" (column C) and "NTW" (column D) columns are removed entirely,
# shifting the remaining columns (KL, HK, NTN) left so the result is:
# year, NTS, KL, HK, NTN (A:E)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D (NTW) first, then column C (NTE), so references stay valid.
$ws.Range("D1:D18").EntireColumn.Delete() | Out-Null
$ws.Range("C1:C18").EntireColumn.Delete() | Out-Null

# Update selection to match the recorded post-edit UI state.
$ws.Range("G9").Select() | Out-Null
